$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.326.53"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "1.710.94"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.54"
$ws.Range("E5").Value = "  -0.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5292"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.004"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.06676"
$ws.Range("E8").Value = "  +1.36%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2661"
$ws.Range("E9").Value = "  -0.32%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.83"
$ws.Range("E10").Value = "  -3.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07701"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.501"
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("D13").Value = "1.946.50"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "1.700.91"
$ws.Range("E14").Value = "  -1.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5850"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "0.0₅8214"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.03"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "27.355.01"
$ws.Range("E18").Value = "  -0.66%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "222.63"
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.005"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.631"
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("E23").Value = "  -1.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.48"
$ws.Range("E25").Value = "  -1.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.693"
$ws.Range("E26").Value = "  -2.44%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1206"
$ws.Range("E27").Value = "  -2.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.241"
$ws.Range("E28").Value = "  -2.26%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.24"
$ws.Range("E29").Value = "  -2.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05325"
$ws.Range("E30").Value = "  -3.87%  "
$ws.Range("E31").Value = "  -0.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.460"
$ws.Range("E32").Value = "  -2.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.424"
$ws.Range("E33").Value = "  -0.54%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.636"
$ws.Range("E34").Value = "  -1.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.871"
$ws.Range("E35").Value = "  +0.41%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9528"
$ws.Range("E36").Value = "  -1.37%  "
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5853"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("D39").Value = "1.145.03"
$ws.Range("E39").Value = "  +8.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01635"
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.790"
$ws.Range("E41").Value = "  -2.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.004"
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8399"
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.32"
$ws.Range("D45").Value = "1.853.59"
$ws.Range("E45").Value = "  -0.70%  "
$ws.Range("E46").Value = "  -0.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.67"
$ws.Range("E47").Value = "  -2.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4551"
$ws.Range("E48").Value = "  +2.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.003"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.111"
$ws.Range("E50").Value = "  -1.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05222"
$ws.Range("E51").Value = "  -0.41%  "
